# Add new worksheet "medbert2" (sheetId 6) with K-fold balanced-data results,
# positioned after the existing "medbert1" sheet (last tab).

$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "medbert2"

# Move the freshly added sheet (currently first) to the end of the tab strip,
# i.e. right after "medbert1". Re-fetch the worksheet by name afterwards since
# the old reference tracks tab *position*, not identity, across the move.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)
$ws = $wb.Worksheets.Item("medbert2")

# ---- Header row (row 1), columns B..N ----
$headers = @("Fold", "Version", "Epoch", "Recall", "Precision", "Accuracy", "Fbeta", "Best Recall", "Best Precision", "Best Threshold", "False Neg(0.5)", "False Pos(0.5)", "Val loss")
$headerCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---- Data rows (2..5) ----
# Column A ("fold index") is bold/centered/bordered like the header.
$foldIndex = @(0, 1, 2, 3)
for ($r = 0; $r -lt $foldIndex.Length; $r++) {
    $row = $r + 2
    $cell = $ws.Range("A" + $row)
    $cell.Value = $foldIndex[$r]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# B: Fold name
$fold = @("fold_0", "fold_1", "fold_2", "fold_3")
# C: Version (run timestamp) - plain text
$version = @("19.09_09.08", "19.09_09.15", "19.09_09.23", "19.09_09.30")
# D: Epoch - integer
$epoch = @(6, 8, 7, 7)
# E: Recall - stored as text
$recall = @("0.6363636", "0.7922078", "0.8701299", "0.7948718")
# F: Precision - stored as text
$precision = @("0.35766423", "0.36746988", "0.34358975", "0.3668639")
# G: Accuracy - stored as text
$accuracy = @("0.8987784", "0.8944154", "0.87958115", "0.89267015")
# H: Fbeta - stored as text
$fbeta = @("0.5505618", "0.6434599", "0.666004", "0.64449066")
# I: Best Recall - number
$bestRecall = @(0, 0, 0, 0)
# J: Best Precision - number
$bestPrecision = @(0, 0, 0, 0)
# K: Best Threshold - number
$bestThreshold = @(0.5, 0.5, 0.5, 0.5)
# L: False Neg(0.5) - integer
$falseNeg = @(28, 16, 10, 16)
# M: False Pos(0.5) - integer
$falsePos = @(88, 105, 128, 107)
# N: Val loss - number
$valLoss = @(0.7052134337524573, 0.3302207328379154, 0.7126566295822462, 0.3642424587160349)

for ($r = 0; $r -lt 4; $r++) {
    $row = $r + 2

    $ws.Range("B" + $row).Value = $fold[$r]

    $ws.Range("C" + $row).NumberFormat = "@"
    $ws.Range("C" + $row).Value = $version[$r]

    $ws.Range("D" + $row).Value = $epoch[$r]

    $ws.Range("E" + $row).NumberFormat = "@"
    $ws.Range("E" + $row).Value = $recall[$r]

    $ws.Range("F" + $row).NumberFormat = "@"
    $ws.Range("F" + $row).Value = $precision[$r]

    $ws.Range("G" + $row).NumberFormat = "@"
    $ws.Range("G" + $row).Value = $accuracy[$r]

    $ws.Range("H" + $row).NumberFormat = "@"
    $ws.Range("H" + $row).Value = $fbeta[$r]

    $ws.Range("I" + $row).Value = $bestRecall[$r]
    $ws.Range("J" + $row).Value = $bestPrecision[$r]
    $ws.Range("K" + $row).Value = $bestThreshold[$r]
    $ws.Range("L" + $row).Value = $falseNeg[$r]
    $ws.Range("M" + $row).Value = $falsePos[$r]
    $ws.Range("N" + $row).Value = $valLoss[$r]
}
